# Adds columns I (I0) and J (IF): header labels + per-row integer values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells, matching style of existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Per-row data: row number, I value, J value
$data = @(
    @(2, 9, 9),
    @(3, 5, 6),
    @(4, 8, 9),
    @(5, 6, 7),
    @(6, 8, 8),
    @(7, 4, 6),
    @(8, 4, 5),
    @(9, 6, 6),
    @(10, 5, 6),
    @(11, 6, 8),
    @(12, 6, 6),
    @(13, 7, 7),
    @(14, 7, 7),
    @(15, 7, 7),
    @(16, 8, 8),
    @(17, 6, 6),
    @(18, 6, 7),
    @(19, 6, 7),
    @(20, 7, 7),
    @(21, 8, 8),
    @(22, 7, 8),
    @(23, 8, 9),
    @(24, 6, 7),
    @(25, 8, 8),
    @(26, 8, 8),
    @(27, 7, 8),
    @(28, 10, 10),
    @(29, 8, 8),
    @(30, 7, 7),
    @(31, 9, 9),
    @(32, 7, 7),
    @(33, 7, 7),
    @(34, 6, 7),
    @(35, 6, 6),
    @(36, 6, 7),
    @(37, 7, 7),
    @(38, 8, 8),
    @(39, 9, 9),
    @(40, 7, 7),
    @(41, 9, 9),
    @(42, 9, 9),
    @(43, 5, 5),
    @(44, 3, 3),
    @(45, 6, 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
